$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the new BOM row (row 22) for the terminal block datasheet entry.
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "J"
$ws.Range("C22").Value = "Anschlussklemme"
$ws.Range("D22").Value = "J2"
$ws.Range("E22").Value = "AKL 057-02"

# Grow the worksheet table (Tabelle1) so it covers the new row.
$table = $ws.ListObjects.Item("Tabelle1")
$table.Resize($ws.Range("A1:G22"))

# Mirror the final selection state recorded for the sheet.
$ws.Range("I30:J30").Select()
